$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 122, pushing the existing row 122
# (and everything below it) down by one row.
$ws.Rows.Item(122).Insert()

# Populate the newly inserted row 122 with a new weekly price record.
$ws.Cells.Item(122, 1).Value  = 4
$ws.Cells.Item(122, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(122, 3).Value  = "Los Lagos"
$ws.Cells.Item(122, 4).Value  = 44705
$ws.Cells.Item(122, 5).Value  = 10
$ws.Cells.Item(122, 6).Value  = 100112017
$ws.Cells.Item(122, 7).Value  = "Apio"
$ws.Cells.Item(122, 8).Value  = "Americana (o)"
$ws.Cells.Item(122, 9).Value  = "Primera"
$ws.Cells.Item(122, 10).Value = 45
$ws.Cells.Item(122, 11).Value = 11000
$ws.Cells.Item(122, 12).Value = 11000
$ws.Cells.Item(122, 13).Value = 11000
$ws.Cells.Item(122, 14).Value = "$/docena de matas"
$ws.Cells.Item(122, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(122, 16).Value = 1833
$ws.Cells.Item(122, 17).Value = 6
$ws.Cells.Item(122, 18).Value = "Hortaliza"
